$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.170827865600586
$ws.Range("B1").Value = 2.472141742706299
$ws.Range("C1").Value = 6.615234851837158
$ws.Range("D1").Value = 2.057464838027954
$ws.Range("E1").Value = 1.206183075904846
